$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Today`n68°`n/59°`nMostly Sunny`n2%`nNW 12 mph"
$ws.Range("A2").Value = "Thu 11`n83°`n/61°`nSunny`n1%`nNNW 9 mph"
$ws.Range("A3").Value = "Fri 12`n85°`n/61°`nSunny`n0%`nNNW 9 mph"
$ws.Range("A4").Value = "Sat 13`n81°`n/61°`nSunny`n1%`nWNW 8 mph"
$ws.Range("A5").Value = "Sun 14`n79°`n/59°`nSunny`n1%`nNW 8 mph"
$ws.Range("A6").Value = "Mon 15`n76°`n/58°`nMostly Sunny`n2%`nWNW 8 mph"
$ws.Range("A7").Value = "Tue 16`n69°`n/58°`nMostly Cloudy`n4%`nWSW 8 mph"
$ws.Range("A8").Value = "Wed 17`n66°`n/56°`nPartly Cloudy`n5%`nSSW 9 mph"
$ws.Range("A9").Value = "Thu 18`n67°`n/55°`nPartly Cloudy`n4%`nWSW 8 mph"
$ws.Range("A10").Value = "Fri 19`n69°`n/56°`nPartly Cloudy`n3%`nNNW 9 mph"
$ws.Range("A11").Value = "Sat 20`n68°`n/57°`nPartly Cloudy`n1%`nNNW 9 mph"
$ws.Range("A12").Value = "Sun 21`n69°`n/58°`nPartly Cloudy`n5%`nNNW 9 mph"
$ws.Range("A13").Value = "Mon 22`n69°`n/57°`nPartly Cloudy`n1%`nNW 8 mph"
$ws.Range("A14").Value = "Tue 23`n68°`n/57°`nPartly Cloudy`n1%`nWNW 8 mph"
$ws.Range("A15").Value = "Wed 24`n67°`n/57°`nMostly Sunny`n1%`nWNW 9 mph"
